$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: "auto" command - run autosequence, takes 1 arg (autosequence commands (list))
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "auto"
$ws.Range("C9").Value = "run autosequence"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "autosequence commands (list)"

# New row 10: "abort_auto" command - abort autosequence, takes 0 args
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "abort_auto"
$ws.Range("C10").Value = "abort autosequence"
$ws.Range("D10").Value = 0

# Update selection to E10 (matches the diff's new active cell)
$ws.Range("E10").Select()
